# #9 option to replace tab with replacement string
# Adds a demonstration row/cell showing a value that contains a literal
# tab character, stored as a new shared string and placed in D2 on the
# "Characters" sheet (between the existing Email and Header-with-linebreak
# sample cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value2 = "Value with`ttab"

# Give the new cell its own (distinct) cell style, same as a fresh/default
# format, so it gets its own entry in the stylesheet rather than continuing
# to share the style used by neighbouring cells.
$cell.ClearFormats() | Out-Null

# Match the author's resulting selection (the newly-edited cell).
$cell.Select() | Out-Null
